$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in row 1 for columns P and Q
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Match formatting of the existing header cells (bold font, border, centered)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update data rows 2 through 25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # column I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # column K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # column M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # column O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # column P: new
    $ws.Cells.Item($r, 17).Value = 2  # column Q: new
}
